# Early_preterm model selection.xlsx -- apply the "Add files via upload" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Main results grid C3:G12 (model picked per fold/run) -- values changed
# ---------------------------------------------------------------------------
$ws.Range("C3").Value  = "BernoulliNB"
$ws.Range("D3").Value  = "CalibratedClassifierCV"
$ws.Range("E3").Value  = "SVC"
$ws.Range("F3").Value  = "SVC"
$ws.Range("G3").Value  = "SVC"

$ws.Range("C4").Value  = "NearestCentroid"
$ws.Range("D4").Value  = "SVC"
$ws.Range("E4").Value  = "QuadraticDiscriminantAnalysis"
$ws.Range("F4").Value  = "SGDClassifier"
$ws.Range("G4").Value  = "QuadraticDiscriminantAnalysis"

$ws.Range("C5").Value  = "ExtraTreeClassifier"
$ws.Range("D5").Value  = "PassiveAggressiveClassifier"
$ws.Range("E5").Value  = "PassiveAggressiveClassifier"
$ws.Range("F5").Value  = "Perceptron"
$ws.Range("G5").Value  = "LinearSVC"

$ws.Range("C6").Value  = "RidgeClassifierCV"
$ws.Range("D6").Value  = "LogisticRegression"
$ws.Range("E6").Value  = "CalibratedClassifierCV"
$ws.Range("F6").Value  = "PassiveAggressiveClassifier"
$ws.Range("G6").Value  = "PassiveAggressiveClassifier"

$ws.Range("C7").Value  = "ExtraTreesClassifier"
$ws.Range("D7").Value  = "LinearSVC"
$ws.Range("E7").Value  = "LinearSVC"
$ws.Range("F7").Value  = "GaussianNB"
$ws.Range("G7").Value  = "CalibratedClassifierCV"

$ws.Range("C8").Value  = "PassiveAggressiveClassifier"
$ws.Range("D8").Value  = "NearestCentroid"
$ws.Range("E8").Value  = "Perceptron"
$ws.Range("F8").Value  = "QuadraticDiscriminantAnalysis"
$ws.Range("G8").Value  = "SGDClassifier"

$ws.Range("C9").Value  = "LogisticRegression"
$ws.Range("D9").Value  = "RidgeClassifierCV"
$ws.Range("E9").Value  = "ExtraTreeClassifier"
$ws.Range("F9").Value  = "AdaBoostClassifier"
$ws.Range("G9").Value  = "Perceptron"

$ws.Range("C10").Value = "LinearSVC"
$ws.Range("D10").Value = "ExtraTreeClassifier"
$ws.Range("E10").Value = "BernoulliNB"
$ws.Range("F10").Value = "ExtraTreeClassifier"
$ws.Range("G10").Value = "LogisticRegression"

$ws.Range("C11").Value = "RidgeClassifier"
$ws.Range("D11").Value = "RidgeClassifier"
$ws.Range("E11").Value = "BaggingClassifier"
$ws.Range("F11").Value = "RidgeClassifierCV"
$ws.Range("G11").Value = "GaussianNB"

$ws.Range("C12").Value = "CalibratedClassifierCV"
$ws.Range("D12").Value = "QuadraticDiscriminantAnalysis"
$ws.Range("E12").Value = "GaussianNB"
$ws.Range("F12").Value = "LinearSVC"
$ws.Range("G12").Value = "LinearDiscriminantAnalysis"

# ---------------------------------------------------------------------------
# 2. Ensemble accuracy summary, row 13
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "0.94, 0.84, 0.83"
$ws.Range("D13").Value = "0.85, 0.85, 0.85"
$ws.Range("E13").Value = "0.85, 0.81, 0.74"
$ws.Range("F13").Value = "0.79, 0.74, 0.68"
$ws.Range("G13").Value = "0.74, 0.66, 0.65"

# ---------------------------------------------------------------------------
# 3. Frequency-count table, column B (labels re-shuffled); COUNTIF formulas
#    in column C / F recalc automatically once the grid + labels change.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "BernoulliNB"
$ws.Range("B18").Value = "NearestCentroid"
$ws.Range("B19").Value = "ExtraTreeClassifier"
$ws.Range("B20").Value = "RidgeClassifierCV"
$ws.Range("B21").Value = "ExtraTreesClassifier"
$ws.Range("B22").Value = "PassiveAggressiveClassifier"
$ws.Range("B23").Value = "LogisticRegression"
$ws.Range("B24").Value = "LinearSVC"
$ws.Range("B25").Value = "RidgeClassifier"
$ws.Range("B26").Value = "CalibratedClassifierCV"
$ws.Range("B27").Value = "SVC"
$ws.Range("B28").Value = "QuadraticDiscriminantAnalysis"
$ws.Range("B29").Value = "Perceptron"
$ws.Range("B30").Value = "BaggingClassifier"
$ws.Range("B31").Value = "GaussianNB"
$ws.Range("B32").Value = "SGDClassifier"
$ws.Range("B33").Value = "AdaBoostClassifier"

# New row 34: one more classifier label + its COUNTIF formula
$ws.Range("B34").Value = "LinearDiscriminantAnalysis"
$ws.Range("C34").Formula = "=COUNTIF(`$C`$3:`$G`$12, B34)"

# ---------------------------------------------------------------------------
# 4. New blank styled cells that appeared to the left of / inside the grid
#    (A6, A9, A11, A12, J11) -- same bold-Arial-9pt look as the grid cells.
# ---------------------------------------------------------------------------
foreach ($addr in "A6","A9","A11","A12","J11") {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 9
    $rng.Font.Color = 0
}

# ---------------------------------------------------------------------------
# 5. Column A width (new column inserted to the left of the table)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 22.43

# ---------------------------------------------------------------------------
# 6. Sheet view: drop the old scrolled/selected cell, select J11 instead
# ---------------------------------------------------------------------------
$ws.Range("J11").Select()

Write-Host "done"
